# Update sheet data cells with newly recomputed TPM-based values.
# Source: commit "update scripts wuth new tpm" - recalculated ligand/receptor
# average & total expression values (TPM-normalized) and all values that
# derive from them (specificity scores, edge weights, edge specificities).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.044933333333333
$ws.Range("H2").Value = 18.1348
$ws.Range("I2").Value = 0.9708761253868625
$ws.Range("J2").Value = 0.9708761253868624
$ws.Range("M2").Value = 8.489835333333334
$ws.Range("N2").Value = 25.469506
$ws.Range("O2").Value = 0.2075776945087381
$ws.Range("P2").Value = 0.2075776945087381
$ws.Range("Q2").Value = 51.32048860097778
$ws.Range("R2").Value = 461.8843974087999
$ws.Range("S2").Value = 0.2015322277613815
$ws.Range("T2").Value = 0.2015322277613814
$ws.Range("G3").Value = 6.044933333333333
$ws.Range("H3").Value = 18.1348
$ws.Range("I3").Value = 0.9708761253868625
$ws.Range("J3").Value = 0.9708761253868624
$ws.Range("O3").Value = 0.3214784855238645
$ws.Range("P3").Value = 0.3214784855238645
$ws.Range("Q3").Value = 79.48076015986666
$ws.Range("R3").Value = 715.3268414387999
$ws.Range("S3").Value = 0.3121157864206461
$ws.Range("T3").Value = 0.3121157864206461
$ws.Range("G4").Value = 6.044933333333333
$ws.Range("H4").Value = 18.1348
$ws.Range("I4").Value = 0.9708761253868625
$ws.Range("J4").Value = 0.9708761253868624
$ws.Range("M4").Value = 5.630791333333334
$ws.Range("N4").Value = 16.892374
$ws.Range("O4").Value = 0.1376736576555254
$ws.Range("P4").Value = 0.1376736576555254
$ws.Range("Q4").Value = 34.03775822391111
$ws.Range("R4").Value = 306.3398240152
$ws.Range("S4").Value = 0.1336640673124339
$ws.Range("T4").Value = 0.1336640673124339
$ws.Range("G5").Value = 6.044933333333333
$ws.Range("H5").Value = 18.1348
$ws.Range("I5").Value = 0.9708761253868625
$ws.Range("J5").Value = 0.9708761253868624
$ws.Range("M5").Value = 6.738585333333333
$ws.Range("N5").Value = 20.215756
$ws.Range("O5").Value = 0.1647593802263456
$ws.Range("P5").Value = 0.1647593802263456
$ws.Range("Q5").Value = 40.73429910097777
$ws.Range("R5").Value = 366.6086919087999
$ws.Range("S5").Value = 0.1599609486952952
$ws.Range("T5").Value = 0.1599609486952952
$ws.Range("G6").Value = 6.044933333333333
$ws.Range("H6").Value = 18.1348
$ws.Range("I6").Value = 0.9708761253868625
$ws.Range("J6").Value = 0.9708761253868624
$ws.Range("M6").Value = 0.9376886666666667
$ws.Range("N6").Value = 2.813066
$ws.Range("O6").Value = 0.02292662271427322
$ws.Range("P6").Value = 0.02292662271427321
$ws.Range("Q6").Value = 5.668265477422223
$ws.Range("R6").Value = 51.0143892968
$ws.Range("S6").Value = 0.02225891062904001
$ws.Range("T6").Value = 0.02225891062904001
$ws.Range("G7").Value = 6.044933333333333
$ws.Range("H7").Value = 18.1348
$ws.Range("I7").Value = 0.9708761253868625
$ws.Range("J7").Value = 0.9708761253868624
$ws.Range("M7").Value = 5.954327333333333
$ws.Range("N7").Value = 17.862982
$ws.Range("O7").Value = 0.1455841593712531
$ws.Range("P7").Value = 0.1455841593712531
$ws.Range("Q7").Value = 35.99351177484444
$ws.Range("R7").Value = 323.9416059736
$ws.Range("S7").Value = 0.1413441845680657
$ws.Range("T7").Value = 0.1413441845680657
$ws.Range("I8").Value = 0.01821359071319307
$ws.Range("J8").Value = 0.01821359071319307
$ws.Range("M8").Value = 8.489835333333334
$ws.Range("N8").Value = 25.469506
$ws.Range("O8").Value = 0.2075776945087381
$ws.Range("P8").Value = 0.2075776945087381
$ws.Range("Q8").Value = 0.962769966360889
$ws.Range("R8").Value = 8.664929697248001
$ws.Range("S8").Value = 0.003780735168970381
$ws.Range("T8").Value = 0.003780735168970381
$ws.Range("I9").Value = 0.01821359071319307
$ws.Range("J9").Value = 0.01821359071319307
$ws.Range("O9").Value = 0.3214784855238645
$ws.Range("P9").Value = 0.3214784855238645
$ws.Range("S9").Value = 0.005855277558428832
$ws.Range("T9").Value = 0.005855277558428832
$ws.Range("I10").Value = 0.01821359071319307
$ws.Range("J10").Value = 0.01821359071319307
$ws.Range("M10").Value = 5.630791333333334
$ws.Range("N10").Value = 16.892374
$ws.Range("O10").Value = 0.1376736576555254
$ws.Range("P10").Value = 0.1376736576555254
$ws.Range("Q10").Value = 0.6385467526435555
$ws.Range("R10").Value = 5.746920773792001
$ws.Range("S10").Value = 0.002507531652526001
$ws.Range("T10").Value = 0.002507531652526
$ws.Range("I11").Value = 0.01821359071319307
$ws.Range("J11").Value = 0.01821359071319307
$ws.Range("M11").Value = 6.738585333333333
$ws.Range("N11").Value = 20.215756
$ws.Range("O11").Value = 0.1647593802263456
$ws.Range("P11").Value = 0.1647593802263456
$ws.Range("Q11").Value = 0.7641735463608889
$ws.Range("R11").Value = 6.877561917248
$ws.Range("S11").Value = 0.003000859917602014
$ws.Range("T11").Value = 0.003000859917602014
$ws.Range("I12").Value = 0.01821359071319307
$ws.Range("J12").Value = 0.01821359071319307
$ws.Range("M12").Value = 0.9376886666666667
$ws.Range("N12").Value = 2.813066
$ws.Range("O12").Value = 0.02292662271427322
$ws.Range("P12").Value = 0.02292662271427321
$ws.Range("Q12").Value = 0.1063363953031111
$ws.Range("R12").Value = 0.9570275577280001
$ws.Range("S12").Value = 0.000417576122553568
$ws.Range("T12").Value = 0.0004175761225535679
$ws.Range("I13").Value = 0.01821359071319307
$ws.Range("J13").Value = 0.01821359071319307
$ws.Range("M13").Value = 5.954327333333333
$ws.Range("N13").Value = 17.862982
$ws.Range("O13").Value = 0.1455841593712531
$ws.Range("P13").Value = 0.1455841593712531
$ws.Range("Q13").Value = 0.6752365978062221
$ws.Range("R13").Value = 6.077129380255999
$ws.Range("S13").Value = 0.002651610293112276
$ws.Range("T13").Value = 0.002651610293112276
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.06793033333333333
$ws.Range("H14").Value = 0.203791
$ws.Range("I14").Value = 0.01091028389994453
$ws.Range("J14").Value = 0.01091028389994453
$ws.Range("M14").Value = 8.489835333333334
$ws.Range("N14").Value = 25.469506
$ws.Range("O14").Value = 0.2075776945087381
$ws.Range("P14").Value = 0.2075776945087381
$ws.Range("Q14").Value = 0.5767173441384444
$ws.Range("R14").Value = 5.190456097246
$ws.Range("S14").Value = 0.00226473157838629
$ws.Range("T14").Value = 0.00226473157838629
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.06793033333333333
$ws.Range("H15").Value = 0.203791
$ws.Range("I15").Value = 0.01091028389994453
$ws.Range("J15").Value = 0.01091028389994453
$ws.Range("O15").Value = 0.3214784855238645
$ws.Range("P15").Value = 0.3214784855238645
$ws.Range("Q15").Value = 0.8931702358856666
$ws.Range("R15").Value = 8.038532122971001
$ws.Range("S15").Value = 0.00350742154478957
$ws.Range("T15").Value = 0.00350742154478957
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.06793033333333333
$ws.Range("H16").Value = 0.203791
$ws.Range("I16").Value = 0.01091028389994453
$ws.Range("J16").Value = 0.01091028389994453
$ws.Range("M16").Value = 5.630791333333334
$ws.Range("N16").Value = 16.892374
$ws.Range("O16").Value = 0.1376736576555254
$ws.Range("P16").Value = 0.1376736576555254
$ws.Range("Q16").Value = 0.3825015322037777
$ws.Range("R16").Value = 3.442513789834
$ws.Range("S16").Value = 0.001502058690565555
$ws.Range("T16").Value = 0.001502058690565554
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.06793033333333333
$ws.Range("H17").Value = 0.203791
$ws.Range("I17").Value = 0.01091028389994453
$ws.Range("J17").Value = 0.01091028389994453
$ws.Range("M17").Value = 6.738585333333333
$ws.Range("N17").Value = 20.215756
$ws.Range("O17").Value = 0.1647593802263456
$ws.Range("P17").Value = 0.1647593802263456
$ws.Range("Q17").Value = 0.4577543478884444
$ws.Range("R17").Value = 4.119789130996
$ws.Range("S17").Value = 0.001797571613448337
$ws.Range("T17").Value = 0.001797571613448337
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.06793033333333333
$ws.Range("H18").Value = 0.203791
$ws.Range("I18").Value = 0.01091028389994453
$ws.Range("J18").Value = 0.01091028389994453
$ws.Range("M18").Value = 0.9376886666666667
$ws.Range("N18").Value = 2.813066
$ws.Range("O18").Value = 0.02292662271427322
$ws.Range("P18").Value = 0.02292662271427321
$ws.Range("Q18").Value = 0.06369750368955555
$ws.Range("R18").Value = 0.5732775332060001
$ws.Range("S18").Value = 0.0002501359626796377
$ws.Range("T18").Value = 0.0002501359626796376
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.06793033333333333
$ws.Range("H19").Value = 0.203791
$ws.Range("I19").Value = 0.01091028389994453
$ws.Range("J19").Value = 0.01091028389994453
$ws.Range("M19").Value = 5.954327333333333
$ws.Range("N19").Value = 17.862982
$ws.Range("O19").Value = 0.1455841593712531
$ws.Range("P19").Value = 0.1455841593712531
$ws.Range("Q19").Value = 0.4044794405291111
$ws.Range("R19").Value = 3.640314964762
$ws.Range("S19").Value = 0.001588364510075142
$ws.Range("T19").Value = 0.001588364510075142
